$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Objectives:/EN text) is untouched by this edit, so its cells keep the
# original style indices and make safe templates to copy-format from; this avoids
# synthesising brand-new style entries in styles.xml for re-used cells.
$styleA = $ws.Range("A11")   # bold label style (s=1)
$styleB = $ws.Range("B11")   # wrap text style   (s=2)
$styleC = $ws.Range("C11")   # red wrap style    (s=3)

# --- Row 10 ---
$ws.Range("A10:C10").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 'Objetivos:'
$styleB.Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = 'Proporcionar ao graduando conhecimentos da teoria básica dos conceitos de transferência de massa com posterior aplicação aos balanços de massa visando obtenção, para os diversos processos físicos e químicos, em particularidade os sistemas estagnados e convectivos, conhecimento do fluxo de transferência de massa, do perfil de concentração, das resistências que prediz o transporte entre as fases.'
$styleC.Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = 'Proporcionar ao graduando conhecimentos da teoria básica dos conceitos de transferência de massa com posterior aplicação aos balanços de massa visando obtenção, para os diversos processos físicos e químicos, em particularidade os sistemas estagnados e convectivos, conhecimento do fluxo de transferência de massa, do perfil de concentração, das resistências que prediz o transporte entre as fases.'
$ws.Rows.Item(10).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 13 ---
$ws.Range("A13:C13").Clear()
$styleB.Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = '5840841 - Gilberto Garcia Cortez'
$styleC.Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = '5840841 - Gilberto Garcia Cortez'
$ws.Rows.Item(13).AutoFit()
$excel.CutCopyMode = $false

# --- Row 14 ---
$ws.Range("A14:C14").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = 'Programa resumido:'
$styleB.Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = '1- Introdução:2- Coeficiente de difusão:3- Concentrações, velocidade e fluxos:4 -Equações da continuidade em transferência de massa:5- Difusão em regime permanente sem reação química:6- Difusão com reação química:7- Transferência de massa entre fases.'
$styleC.Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value = '1- Introdução:2- Coeficiente de difusão:3- Concentrações, velocidade e fluxos:4 -Equações da continuidade em transferência de massa:5- Difusão em regime permanente sem reação química:6- Difusão com reação química:7- Transferência de massa entre fases.'
$ws.Rows.Item(14).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 15 ---
$ws.Range("A15:C15").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = 'Short syllabus:'
$styleB.Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = '1 - Introduction: 2 - Diffusion coefficient: 3 - Concentrations, and flow rate: 4 - Equation of continuity for mass transfer: 5 - Diffusion in continuous operation without chemical reaction: 6 - Diffusion with chemical reaction: 7 - Mass transfer between phases.'
$styleC.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = '1 - Introduction: 2 - Diffusion coefficient: 3 - Concentrations, and flow rate: 4 - Equation of continuity for mass transfer: 5 - Diffusion in continuous operation without chemical reaction: 6 - Diffusion with chemical reaction: 7 - Mass transfer between phases.'
$ws.Rows.Item(15).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 16 ---
$ws.Range("A16:C16").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = 'Programa:'
$styleB.Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = '1- Introdução: Transferência de massa: Definição. Classificação das operações que envolvem transferência de massa. Contribuições à transferência de massa. Tipos de difusão.2- Coeficiente de difusão: Considerações. Difusão em gases: Análise da primeira lei de Fick; O coeficiente de difusão para gases. Estimativa do coeficiente de difusão a partir de um coeficiente de difusão conhecido em alta temperatura e pressão. Coeficiente de difusão de um soluto em uma mistura gasosa estagnada de multicomponentes. Difusão em líquidos. Difusão em sólidos.3- Concentrações, velocidade e fluxos: Concentração. Velocidade. Fluxo. A equação de Stefan – Maxwel.4 - Equações da continuidade em transferência de massa: Equações da continuidade molar de um soluto. Regime transiente sem/com velocidade do meio nula. Meio sem e com reação química.5- Difusão em regime permanente sem reação química: Difusão Unidimensional em regime permanente. Difusão através de filme gasoso inerte e estagnado. Difusão pseudo-estacionária num filme gasoso estagnado. Contradifusão equimolar. Taxa molar em esferas isoladas. Difusão em membranas.6- Difusão em regime permanente com reação química: Difusão em regime permanente com reação química heterogênea na superfície de uma partícula catalítica não porosa. Difusão com reação química heterogênea na superfície de uma partícula não catalítica e não porosa. Difusão intrapartícula com reação química heterogênea. Difusão em regime permanente com reação química homogênea.7- Transferência de massa entre fases: Teoria das duas resistências. Coeficiente individual e global de transferência de massa. Coeficientes globais de transferência de massa. Coeficientes volumétricos de transferência de massa para torres de recheios. Balanço macroscópio de matéria. Operações contínuas.'
$styleC.Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = '1- Introdução: Transferência de massa: Definição. Classificação das operações que envolvem transferência de massa. Contribuições à transferência de massa. Tipos de difusão.2- Coeficiente de difusão: Considerações. Difusão em gases: Análise da primeira lei de Fick; O coeficiente de difusão para gases. Estimativa do coeficiente de difusão a partir de um coeficiente de difusão conhecido em alta temperatura e pressão. Coeficiente de difusão de um soluto em uma mistura gasosa estagnada de multicomponentes. Difusão em líquidos. Difusão em sólidos.3- Concentrações, velocidade e fluxos: Concentração. Velocidade. Fluxo. A equação de Stefan – Maxwel.4 - Equações da continuidade em transferência de massa: Equações da continuidade molar de um soluto. Regime transiente sem/com velocidade do meio nula. Meio sem e com reação química.5- Difusão em regime permanente sem reação química: Difusão Unidimensional em regime permanente. Difusão através de filme gasoso inerte e estagnado. Difusão pseudo-estacionária num filme gasoso estagnado. Contradifusão equimolar. Taxa molar em esferas isoladas. Difusão em membranas.6- Difusão em regime permanente com reação química: Difusão em regime permanente com reação química heterogênea na superfície de uma partícula catalítica não porosa. Difusão com reação química heterogênea na superfície de uma partícula não catalítica e não porosa. Difusão intrapartícula com reação química heterogênea. Difusão em regime permanente com reação química homogênea.7- Transferência de massa entre fases: Teoria das duas resistências. Coeficiente individual e global de transferência de massa. Coeficientes globais de transferência de massa. Coeficientes volumétricos de transferência de massa para torres de recheios. Balanço macroscópio de matéria. Operações contínuas.'
$ws.Rows.Item(16).RowHeight = 120
$excel.CutCopyMode = $false

# --- Row 17 ---
$ws.Range("A17:C17").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 'Syllabus:'
$styleB.Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = '1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick''s first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan – Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.'
$styleC.Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = '1 - Introduction: Mass transfer: Definition. Classification of operations involving mass transfer. Contributions to mass transfer. Types of diffusion. 2 - Diffusion coefficient: Considerations. Diffusion in gases: Analysis of Fick''s first law, the diffusion coefficient for gases. Estimation of the diffusion coefficient from a known diffusion coefficient at high temperature and pressure. Diffusion coefficient of a solute in a multicomponent mixture of stagnant gases. Diffusion in liquids. Diffusion in solids. 3 - Concentrations, and flow rate: Concentration. Speed and Flow. The equation of Stefan – Maxwell. 4 - Equation of continuity for mass transfer: Equations of continuity of a solute molar. Transient without/with speed zero means. Medium with and without chemical reaction. 5 - Diffusion in permanent regime without chemical reaction: One-dimensional diffusion in permanent regime. Diffusion through inert and stagnant gaseous film. Pseudo-stationary diffusion in a stagnant gaseous film. Equimolar contradifusion. Molar rate in isolated beads. Diffusion in membranes. 6 - Diffusion in permanent with chemical reaction: Diffusion in steady state with heterogeneous chemical reaction on the surface of a nonporous catalytic particle. Diffusion with heterogeneous chemical reaction on the surface of a non-catalytic and non-porous particle. Intraparticle diffusion with heterogeneous chemical reaction. Continuous diffusion with homogeneous chemical reaction. 7 - Mass transfer between phases: Theory of the two resistors. Individual and global coefficient of mass transfer. Global mass transfer coefficients. Volumetric coefficients of transfer of mass to towers of fillings. Balance macroscope of matter. Continuous operations.'
$ws.Rows.Item(17).RowHeight = 120
$excel.CutCopyMode = $false

# --- Row 18 ---
$ws.Range("A18:C18").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 'Avaliação:'
$ws.Rows.Item(18).AutoFit()
$excel.CutCopyMode = $false

# --- Row 19 ---
$ws.Range("A19:C19").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = 'Método:'
$styleB.Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$styleC.Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Rows.Item(19).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 20 ---
$ws.Range("A20:C20").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = 'Critério:'
$styleB.Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value = 'A Nota Final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3.'
$styleC.Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 'A Nota Final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3.'
$ws.Rows.Item(20).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 21 ---
$ws.Range("A21:C21").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = 'Norma de recuperação:'
$styleB.Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma Prova Escrita (PE) e a Média de Recuperação (MR) será calculada pela fórmula: MR = (NF + PE)/2.'
$styleC.Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma Prova Escrita (PE) e a Média de Recuperação (MR) será calculada pela fórmula: MR = (NF + PE)/2.'
$ws.Rows.Item(21).RowHeight = 60
$excel.CutCopyMode = $false

# --- Row 22 ---
$ws.Range("A22:C22").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = 'Bibliografia:'
$styleB.Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Value = '1) CREMASCO, M. A. Fundamentos de Transferência de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transferência de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. 2ª ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fenômeno de Transporte: Quantidade de Movimento, Calor e Massa. São Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princípios das Operações Unitárias. 2ª ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2ª ed. Rio de Janeiro: Editora Interciência, 2007.'
$styleC.Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = '1) CREMASCO, M. A. Fundamentos de Transferência de Massa. Campinas: Editora Unicamp; 2008-2009.2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transferência de Calor e Massa. Rio de Janeiro: LTC, 2008.3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. 2ª ed. Rio de Janeiro: LTC, 2004.4) BENNETT, C. O.; MYERS, J. E. Fenômeno de Transporte: Quantidade de Movimento, Calor e Massa. São Paulo: McGrawc- Hill, 1978.5) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 5th ed. Pergamon Press, Oxford, 1996. v.16) FOUST, A. S.; Wenzel, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSON, L. B. Princípios das Operações Unitárias. 2ª ed. Rio de Janeiro: Guanabara Dois, 1982.7) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. New York: McGraw-Hill, 2008.8) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 5th ed. USA: John Wiley & Sons, Inc, 2008.9) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2000.10) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2ª ed. Rio de Janeiro: Editora Interciência, 2007.'
$ws.Rows.Item(22).RowHeight = 120
$excel.CutCopyMode = $false

# --- Row 23 ---
$ws.Range("A23:C23").Clear()
$styleA.Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = 'Requisitos:'
$ws.Rows.Item(23).AutoFit()
$excel.CutCopyMode = $false

# --- Row 24 ---
$ws.Range("A24:C24").Clear()
$styleB.Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = 'LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
'
$styleC.Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value = 'LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30
$excel.CutCopyMode = $false
